# "Don't need to convert with fuhep"
# Column O (Clint, uL/min/10^6 cells) was computed as K*J/120/21.4*1000
# (multiplying Clint-hep by fuhep). Remove the unnecessary fuhep (column J)
# factor so the formula becomes K/120/21.4*1000, fill it down the whole
# table (rows 2-119, including rows where it was previously blank because
# K was blank), and give the column a 2-decimal number format.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("O2:O119").Formula = "=K2/120/21.4*1000"
$ws.Range("O2:O119").NumberFormat = "0.00"

$ws.Range("O2:O119").Select()
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 1
